$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: Taxonsorteringsordning 79592 -> 80258
$ws.Range("B2").Value = 80258

# C2: Valideringsstatus "Ovaliderad" -> cleared (cell removed entirely)
$ws.Range("C2").ClearContents()

# J2, K2, N2: newly materialized empty-text cells (Enhet, Ålder-Stadium, Metod).
# A plain empty-string assignment de-materializes the cell, so we briefly
# force "text" entry mode (leading apostrophe) to create an empty text
# value, then reset the style back to Normal so no stray formatting sticks.
foreach ($addr in @("J2", "K2", "N2")) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

# AC2: append note about microscoping / determination
$ws.Range("AC2").Value = "På gammal solbelyst asp i vägkant. Även spår av asppraktbagge i basen. Mikroskoperad och artbestämd av Fredrik Jonsson."

# AF2: Bestämningsmetod
$ws.Range("AF2").Value = "mikroskoperad"

# AS2: Artbestämd av
$ws.Range("AS2").Value = "Fredrik Jonsson"
